# Hortaliza, Macroferia Regional de Talca - Repollo
# A new weekly price record is inserted at row 131 (pushing every
# subsequent record down by one row, 131->204 instead of 131->203).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 131..203 down to 132..204, leaving a blank row at 131.
$ws.Rows.Item(131).Insert()

# Populate the newly inserted row 131 with the new weekly record.
$ws.Cells.Item(131, 1).Value2 = 5
$ws.Cells.Item(131, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(131, 3).Value2 = "Maule"
$ws.Cells.Item(131, 4).Value2 = 44529
$ws.Cells.Item(131, 5).Value2 = 7
$ws.Cells.Item(131, 6).Value2 = 100112006
$ws.Cells.Item(131, 7).Value2 = "Repollo"
$ws.Cells.Item(131, 8).Value2 = "Crespo record"
$ws.Cells.Item(131, 9).Value2 = "Primera"
$ws.Cells.Item(131, 10).Value2 = 3000
$ws.Cells.Item(131, 11).Value2 = 800
$ws.Cells.Item(131, 12).Value2 = 800
$ws.Cells.Item(131, 13).Value2 = 800
$ws.Cells.Item(131, 14).Value2 = "$/unidad"
$ws.Cells.Item(131, 15).Value2 = "Región del Maule"
$ws.Cells.Item(131, 16).Value2 = 800
$ws.Cells.Item(131, 17).Value2 = 1
$ws.Cells.Item(131, 18).Value2 = "Hortaliza"
